$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("ALC").Range("H98").Value2 = 38462292
$wb.Worksheets.Item("ALC").Range("I98").Value2 = 40000776
$wb.Worksheets.Item("ALC").Range("K98").Value2 = 40000776
$wb.Worksheets.Item("ALC").Range("M98").Value2 = -39999278

$wb.Worksheets.Item("ALC").Range("H106").Value2 = 5189.7
$wb.Worksheets.Item("ALC").Range("I106").Value2 = 3952.5
$wb.Worksheets.Item("ALC").Range("K106").Value2 = 3952.5
$wb.Worksheets.Item("ALC").Range("M106").Value2 = -3321.5

$wb.Worksheets.Item("ALC").Range("H116").Value2 = 7723.5454
$wb.Worksheets.Item("ALC").Range("I116").Value2 = 7244.875
$wb.Worksheets.Item("ALC").Range("K116").Value2 = 7244.875
$wb.Worksheets.Item("ALC").Range("M116").Value2 = -3802.875

$wb.Worksheets.Item("ALC").Range("H122").Value2 = 38462292
$wb.Worksheets.Item("ALC").Range("I122").Value2 = 40000776
$wb.Worksheets.Item("ALC").Range("K122").Value2 = 120002328
$wb.Worksheets.Item("ALC").Range("M122").Value2 = -119999878

$wb.Worksheets.Item("ALC").Range("H131").Value2 = 3031.7
$wb.Worksheets.Item("ALC").Range("I131").Value2 = 1433.3846
$wb.Worksheets.Item("ALC").Range("K131").Value2 = 4300.1538
$wb.Worksheets.Item("ALC").Range("M131").Value2 = 739.8462

$wb.Worksheets.Item("ALC").Range("H132").Value2 = 647.3409
$wb.Worksheets.Item("ALC").Range("I132").Value2 = 639.1395
$wb.Worksheets.Item("ALC").Range("J132").Value2 = 1000
$wb.Worksheets.Item("ALC").Range("K132").Value2 = 1917.4185
$wb.Worksheets.Item("ALC").Range("L132").Value2 = 3000
$wb.Worksheets.Item("ALC").Range("M132").Value2 = 612.5815
$wb.Worksheets.Item("ALC").Range("N132").Value2 = -8060

$wb.Worksheets.Item("ARM").Range("H32").Value2 = 7046664.5
$wb.Worksheets.Item("ARM").Range("I32").Value2 = 8622268
$wb.Worksheets.Item("ARM").Range("J32").Value2 = 17049.691
$wb.Worksheets.Item("ARM").Range("K32").Value2 = 8622268
$wb.Worksheets.Item("ARM").Range("L32").Value2 = 17049.691
$wb.Worksheets.Item("ARM").Range("M32").Value2 = -8621981
$wb.Worksheets.Item("ARM").Range("N32").Value2 = -17623.691

$wb.Worksheets.Item("ARM").Range("H61").Value2 = 13640763
$wb.Worksheets.Item("ARM").Range("I61").Value2 = 10420783
$wb.Worksheets.Item("ARM").Range("J61").Value2 = 35720628
$wb.Worksheets.Item("ARM").Range("K61").Value2 = 10420783
$wb.Worksheets.Item("ARM").Range("L61").Value2 = 35720628
$wb.Worksheets.Item("ARM").Range("M61").Value2 = -10420571
$wb.Worksheets.Item("ARM").Range("N61").Value2 = -35721052

$wb.Worksheets.Item("ARM").Range("H101").Value2 = 117499.5
$wb.Worksheets.Item("ARM").Range("J101").Value2 = 117499.5
$wb.Worksheets.Item("ARM").Range("L101").Value2 = 117499.5
$wb.Worksheets.Item("ARM").Range("N101").Value2 = -123989.5

$wb.Worksheets.Item("ARM").Range("H132").Value2 = 3436.6743
$wb.Worksheets.Item("ARM").Range("I132").Value2 = 1507.3103
$wb.Worksheets.Item("ARM").Range("K132").Value2 = 4521.9309
$wb.Worksheets.Item("ARM").Range("M132").Value2 = -1991.9309

$wb.Worksheets.Item("ARM").Range("H136").Value2 = 13640763
$wb.Worksheets.Item("ARM").Range("I136").Value2 = 10420783
$wb.Worksheets.Item("ARM").Range("J136").Value2 = 35720628
$wb.Worksheets.Item("ARM").Range("K136").Value2 = 31262349
$wb.Worksheets.Item("ARM").Range("L136").Value2 = 107161884
$wb.Worksheets.Item("ARM").Range("M136").Value2 = -31259799
$wb.Worksheets.Item("ARM").Range("N136").Value2 = -107166984

$wb.Worksheets.Item("ARM").Range("H138").Value2 = 174559.58
$wb.Worksheets.Item("ARM").Range("J138").Value2 = 181152.83
$wb.Worksheets.Item("ARM").Range("L138").Value2 = 181152.83
$wb.Worksheets.Item("ARM").Range("N138").Value2 = -191432.83

$wb.Worksheets.Item("BSM").Range("H99").Value2 = 2387.1538
$wb.Worksheets.Item("BSM").Range("I99").Value2 = 1655.5294
$wb.Worksheets.Item("BSM").Range("K99").Value2 = 1655.5294
$wb.Worksheets.Item("BSM").Range("M99").Value2 = -157.5293999999999

$wb.Worksheets.Item("BSM").Range("H134").Value2 = 325195.16
$wb.Worksheets.Item("BSM").Range("I134").Value2 = 2645.44
$wb.Worksheets.Item("BSM").Range("J134").Value2 = 1669152.4
$wb.Worksheets.Item("BSM").Range("K134").Value2 = 7936.32
$wb.Worksheets.Item("BSM").Range("L134").Value2 = 5007457.199999999
$wb.Worksheets.Item("BSM").Range("M134").Value2 = -5401.32
$wb.Worksheets.Item("BSM").Range("N134").Value2 = -5012527.199999999

$wb.Worksheets.Item("CRP").Range("H16").Value2 = 1160.2
$wb.Worksheets.Item("CRP").Range("I16").Value2 = 830
$wb.Worksheets.Item("CRP").Range("J16").Value2 = 1490.4
$wb.Worksheets.Item("CRP").Range("K16").Value2 = 830
$wb.Worksheets.Item("CRP").Range("L16").Value2 = 1490.4
$wb.Worksheets.Item("CRP").Range("M16").Value2 = -543
$wb.Worksheets.Item("CRP").Range("N16").Value2 = -2064.4

$wb.Worksheets.Item("CRP").Range("H31").Value2 = 492769.4
$wb.Worksheets.Item("CRP").Range("I31").Value2 = 6978.2666
$wb.Worksheets.Item("CRP").Range("J31").Value2 = 1302421.4
$wb.Worksheets.Item("CRP").Range("K31").Value2 = 6978.2666
$wb.Worksheets.Item("CRP").Range("L31").Value2 = 1302421.4
$wb.Worksheets.Item("CRP").Range("M31").Value2 = -6683.2666
$wb.Worksheets.Item("CRP").Range("N31").Value2 = -1303011.4

$wb.Worksheets.Item("CRP").Range("H34").Value2 = 492769.4
$wb.Worksheets.Item("CRP").Range("I34").Value2 = 6978.2666
$wb.Worksheets.Item("CRP").Range("J34").Value2 = 1302421.4
$wb.Worksheets.Item("CRP").Range("K34").Value2 = 6978.2666
$wb.Worksheets.Item("CRP").Range("L34").Value2 = 1302421.4
$wb.Worksheets.Item("CRP").Range("M34").Value2 = -6776.2666
$wb.Worksheets.Item("CRP").Range("N34").Value2 = -1302825.4

$wb.Worksheets.Item("CRP").Range("H99").Value2 = 2718.6667
$wb.Worksheets.Item("CRP").Range("J99").Value2 = 3485.375
$wb.Worksheets.Item("CRP").Range("L99").Value2 = 3485.375
$wb.Worksheets.Item("CRP").Range("N99").Value2 = -6481.375

$wb.Worksheets.Item("CRP").Range("H113").Value2 = 1160.2
$wb.Worksheets.Item("CRP").Range("I113").Value2 = 830
$wb.Worksheets.Item("CRP").Range("J113").Value2 = 1490.4
$wb.Worksheets.Item("CRP").Range("K113").Value2 = 830
$wb.Worksheets.Item("CRP").Range("L113").Value2 = 1490.4
$wb.Worksheets.Item("CRP").Range("M113").Value2 = 1340
$wb.Worksheets.Item("CRP").Range("N113").Value2 = -5830.4

$wb.Worksheets.Item("CRP").Range("H122").Value2 = 2420.923
$wb.Worksheets.Item("CRP").Range("I122").Value2 = 2149.125
$wb.Worksheets.Item("CRP").Range("K122").Value2 = 6447.375
$wb.Worksheets.Item("CRP").Range("M122").Value2 = -3997.375

$wb.Worksheets.Item("CRP").Range("H126").Value2 = 2718.6667
$wb.Worksheets.Item("CRP").Range("J126").Value2 = 3485.375
$wb.Worksheets.Item("CRP").Range("L126").Value2 = 10456.125
$wb.Worksheets.Item("CRP").Range("N126").Value2 = -15396.125

$wb.Worksheets.Item("CRP").Range("H132").Value2 = 2193.1538
$wb.Worksheets.Item("CRP").Range("I132").Value2 = 2192.625
$wb.Worksheets.Item("CRP").Range("K132").Value2 = 6577.875
$wb.Worksheets.Item("CRP").Range("M132").Value2 = -4047.875

$wb.Worksheets.Item("CUL").Range("H26").Value2 = 100
$wb.Worksheets.Item("CUL").Range("I26").Value2 = 100
$wb.Worksheets.Item("CUL").Range("K26").Value2 = 300
$wb.Worksheets.Item("CUL").Range("M26").Value2 = -12

$wb.Worksheets.Item("CUL").Range("H131").Value2 = 5431.2554
$wb.Worksheets.Item("CUL").Range("J131").Value2 = 5431.2554
$wb.Worksheets.Item("CUL").Range("L131").Value2 = 16293.7662
$wb.Worksheets.Item("CUL").Range("N131").Value2 = -26373.7662

$wb.Worksheets.Item("GSM").Range("H97").Value2 = 1495.9259
$wb.Worksheets.Item("GSM").Range("I97").Value2 = 1430.65
$wb.Worksheets.Item("GSM").Range("J97").Value2 = 1682.4286
$wb.Worksheets.Item("GSM").Range("K97").Value2 = 1430.65
$wb.Worksheets.Item("GSM").Range("L97").Value2 = 1682.4286
$wb.Worksheets.Item("GSM").Range("M97").Value2 = -934.6500000000001
$wb.Worksheets.Item("GSM").Range("N97").Value2 = -2674.4286

$wb.Worksheets.Item("GSM").Range("H132").Value2 = 14497061
$wb.Worksheets.Item("GSM").Range("I132").Value2 = 20410422
$wb.Worksheets.Item("GSM").Range("K132").Value2 = 61231266
$wb.Worksheets.Item("GSM").Range("M132").Value2 = -61228736

$wb.Worksheets.Item("GSM").Range("H140").Value2 = 80556.5
$wb.Worksheets.Item("GSM").Range("J140").Value2 = 80556.5
$wb.Worksheets.Item("GSM").Range("L140").Value2 = 80556.5
$wb.Worksheets.Item("GSM").Range("N140").Value2 = -90916.5

$wb.Worksheets.Item("LTW").Range("H82").Value2 = 915
$wb.Worksheets.Item("LTW").Range("I82").Value2 = 840
$wb.Worksheets.Item("LTW").Range("J82").Value2 = 990
$wb.Worksheets.Item("LTW").Range("K82").Value2 = 840
$wb.Worksheets.Item("LTW").Range("L82").Value2 = 990
$wb.Worksheets.Item("LTW").Range("M82").Value2 = -479
$wb.Worksheets.Item("LTW").Range("N82").Value2 = -1712

$wb.Worksheets.Item("LTW").Range("H85").Value2 = 915
$wb.Worksheets.Item("LTW").Range("I85").Value2 = 840
$wb.Worksheets.Item("LTW").Range("J85").Value2 = 990
$wb.Worksheets.Item("LTW").Range("K85").Value2 = 840
$wb.Worksheets.Item("LTW").Range("L85").Value2 = 990
$wb.Worksheets.Item("LTW").Range("M85").Value2 = 408
$wb.Worksheets.Item("LTW").Range("N85").Value2 = -3486

$wb.Worksheets.Item("LTW").Range("H93").Value2 = 166667760
$wb.Worksheets.Item("LTW").Range("I93").Value2 = 250000660
$wb.Worksheets.Item("LTW").Range("K93").Value2 = 250000660
$wb.Worksheets.Item("LTW").Range("M93").Value2 = -249999412

$wb.Worksheets.Item("LTW").Range("H100").Value2 = 4138.222
$wb.Worksheets.Item("LTW").Range("I100").Value2 = 2800
$wb.Worksheets.Item("LTW").Range("K100").Value2 = 2800
$wb.Worksheets.Item("LTW").Range("M100").Value2 = -2259

$wb.Worksheets.Item("LTW").Range("H109").Value2 = 100088
$wb.Worksheets.Item("LTW").Range("J109").Value2 = 100088
$wb.Worksheets.Item("LTW").Range("L109").Value2 = 100088
$wb.Worksheets.Item("LTW").Range("N109").Value2 = -102862

$wb.Worksheets.Item("LTW").Range("H123").Value2 = 57975
$wb.Worksheets.Item("LTW").Range("J123").Value2 = 57975
$wb.Worksheets.Item("LTW").Range("L123").Value2 = 57975
$wb.Worksheets.Item("LTW").Range("N123").Value2 = -67775

$wb.Worksheets.Item("LTW").Range("H132").Value2 = 849517.7
$wb.Worksheets.Item("LTW").Range("I132").Value2 = 17655.182
$wb.Worksheets.Item("LTW").Range("J132").Value2 = 10000005
$wb.Worksheets.Item("LTW").Range("K132").Value2 = 52965.546
$wb.Worksheets.Item("LTW").Range("L132").Value2 = 30000015
$wb.Worksheets.Item("LTW").Range("M132").Value2 = -50435.546
$wb.Worksheets.Item("LTW").Range("N132").Value2 = -30005075

$wb.Worksheets.Item("WVR").Range("H52").Value2 = 11788.4
$wb.Worksheets.Item("WVR").Range("I52").Value2 = 12223.75
$wb.Worksheets.Item("WVR").Range("K52").Value2 = 12223.75
$wb.Worksheets.Item("WVR").Range("M52").Value2 = -11997.75

$wb.Worksheets.Item("WVR").Range("H70").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("J70").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("L70").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("N70").Value2 = -100620

$wb.Worksheets.Item("WVR").Range("H73").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("J73").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("L73").Value2 = 99990
$wb.Worksheets.Item("WVR").Range("N73").Value2 = -102174

$wb.Worksheets.Item("WVR").Range("H109").Value2 = 107970
$wb.Worksheets.Item("WVR").Range("J109").Value2 = 107970
$wb.Worksheets.Item("WVR").Range("L109").Value2 = 107970
$wb.Worksheets.Item("WVR").Range("N109").Value2 = -110744

$wb.Worksheets.Item("WVR").Range("H126").Value2 = 979
$wb.Worksheets.Item("WVR").Range("I126").Value2 = 963.2143
$wb.Worksheets.Item("WVR").Range("K126").Value2 = 2889.6429
$wb.Worksheets.Item("WVR").Range("M126").Value2 = -419.6428999999998

$wb.Worksheets.Item("WVR").Range("H132").Value2 = 457826.72
$wb.Worksheets.Item("WVR").Range("I132").Value2 = 3009.3684
$wb.Worksheets.Item("WVR").Range("K132").Value2 = 9028.1052
$wb.Worksheets.Item("WVR").Range("M132").Value2 = -6498.1052

$wb.Worksheets.Item("WVR").Range("H141").Value2 = 68000
$wb.Worksheets.Item("WVR").Range("J141").Value2 = 68000
$wb.Worksheets.Item("WVR").Range("L141").Value2 = 68000
$wb.Worksheets.Item("WVR").Range("N141").Value2 = -78360
